$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-09-01 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-02 Monday", 2) | Out-Null

# Update the division problems in the table, cell by cell (position-based
# to avoid collisions between old/new values that coincide, e.g. 31÷6=)
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "71÷4="
$t.Cell(1,2).Range.Text = "92÷6="
$t.Cell(1,3).Range.Text = "33÷7="
$t.Cell(1,4).Range.Text = "54÷9="
$t.Cell(1,5).Range.Text = "14÷6="
$t.Cell(5,1).Range.Text = "59÷7="
$t.Cell(5,2).Range.Text = "56÷9="
$t.Cell(5,3).Range.Text = "14÷6="
$t.Cell(5,4).Range.Text = "31÷6="
$t.Cell(5,5).Range.Text = "73÷4="
$t.Cell(9,1).Range.Text = "40÷8="
$t.Cell(9,2).Range.Text = "21÷3="
$t.Cell(9,3).Range.Text = "88÷9="
$t.Cell(9,4).Range.Text = "30÷9="
$t.Cell(9,5).Range.Text = "50÷4="
$t.Cell(13,1).Range.Text = "88÷7="
$t.Cell(13,2).Range.Text = "12÷3="
$t.Cell(13,3).Range.Text = "38÷6="
$t.Cell(13,4).Range.Text = "83÷5="
$t.Cell(13,5).Range.Text = "30÷2="
$t.Cell(17,1).Range.Text = "17÷5="
$t.Cell(17,2).Range.Text = "15÷3="
$t.Cell(17,3).Range.Text = "24÷3="
$t.Cell(17,4).Range.Text = "43÷5="
$t.Cell(17,5).Range.Text = "60÷7="
